$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "74.838.74"
$ws.Range("E2").Value = "  +1.27%  "

# Row 3
$ws.Range("D3").Value = "2.816.32"
$ws.Range("E3").Value = "  +7.46%  "

# Row 4
$ws.Range("E4").Value = "  +0.13%  "

# Row 5
$ws.Range("D5").Value = "'187.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.54%  "

# Row 6
$ws.Range("D6").Value = "'591.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.72%  "

# Row 7
$ws.Range("E7").Value = "  +0.09%  "

# Row 8
$ws.Range("E8").Value = "  +3.19%  "

# Row 9
$ws.Range("E9").Value = "  -4.51%  "

# Row 10
$ws.Range("D10").Value = "2.819.66"
$ws.Range("E10").Value = "  +7.72%  "

# Row 11
$ws.Range("D11").Value = "'0.374"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.81%  "

# Row 12
$ws.Range("E12").Value = "  -1.84%  "

# Row 13
$ws.Range("D13").Value = "'4.89"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.51%  "

# Row 15
$ws.Range("D15").Value = "74.810.18"
$ws.Range("E15").Value = "  +1.26%  "

# Row 16
$ws.Range("D16").Value = "'0.0000186"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.29%  "

# Row 17
$ws.Range("D17").Value = "'26.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.80%  "

# Row 18
$ws.Range("D18").Value = "2.814.88"
$ws.Range("E18").Value = "  +6.80%  "

# Row 19
$ws.Range("D19").Value = "'9.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.06%  "

# Row 20
$ws.Range("D20").Value = "'12.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.70%  "

# Row 21
$ws.Range("D21").Value = "'376.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.81%  "

# Row 22
$ws.Range("E22").Value = "  -0.89%  "

# Row 23
$ws.Range("D23").Value = "'4.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.43%  "

# Row 24
$ws.Range("E24").Value = "  +0.08%  "

# Row 25
$ws.Range("D25").Value = "'70.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.96%  "

# Row 26
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "2.947.31"
$ws.Range("E26").Value = "  +6.91%  "

# Row 27
$ws.Range("B27").Value = "NEARProtocol"
$ws.Range("C27").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D27").Value = "'4.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.36%  "

# Row 28
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").Value = "'9.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.35%  "

# Row 29
$ws.Range("D29").Value = "'0.0000103"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +10.15%  "

# Row 30
$ws.Range("D30").Value = "'0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.06%  "

# Row 31
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Value = "'1.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.37%  "

# Row 32
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "'512.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.65%  "

# Row 33
$ws.Range("D33").Value = "'7.66"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.07%  "

# Row 34
$ws.Range("D34").Value = "'1.78"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.08%  "

# Row 35
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.16%  "

# Row 36
$ws.Range("D36").Value = "'164.11"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.95%  "

# Row 37
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "'0.119"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.77%  "

# Row 38
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").Value = "'19.90"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.26%  "

# Row 39
$ws.Range("E39").Value = "  +0.45%  "

# Row 40
$ws.Range("D40").Value = "'185.43"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +14.20%  "

# Row 41
$ws.Range("E41").Value = "  -0.03%  "

# Row 42
$ws.Range("D42").Value = "'0.340"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.67%  "

# Row 43
$ws.Range("D43").Value = "'4.98"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.66%  "

# Row 44
$ws.Range("E44").Value = "  -0.05%  "

# Row 45
$ws.Range("E45").Value = "  +2.58%  "

# Row 46
$ws.Range("D46").Value = "'40.03"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.87%  "

# Row 47
$ws.Range("D47").Value = "'0.0861"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.63%  "

# Row 48
$ws.Range("D48").Value = "'2.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.16%  "

# Row 49
$ws.Range("D49").Value = "'0.574"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +9.71%  "

# Row 50
$ws.Range("E50").Value = "  +3.22%  "

# Row 51
$ws.Range("D51").Value = "'0.636"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +8.89%  "
